# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.028.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.644.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.06%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("E5").Value = "  +0.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "201.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.221"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.57%  "

$ws.Range("E10").Value = "  +0.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000307"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.215.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "678.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +13.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "71.095.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.639.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.99%  "

$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("E21").Value = "  +2.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "105.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.71%  "

$ws.Range("E26").Value = "  -1.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.19%  "

$ws.Range("E33").Value = "  +1.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0869"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.924.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "520.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.12%  "

$ws.Range("E39").Value = "  -5.40%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.53%  "

$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.391"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.137"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0459"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.45%  "

$ws.Range("E45").Value = "  +9.43%  "

$ws.Range("E46").Value = "  +7.29%  "

$ws.Range("E47").Value = "  +1.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.41%  "

$ws.Range("E49").Value = "  -0.39%  "

$ws.Range("E50").Value = "  +2.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.92%  "
